# Open MySQL BVT test class
# Adds three new DML update/delete test case rows (97-99) to Sheet1,
# mirroring the style of the existing updel_09x rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97: updel_096 - 两端范围删除 (range delete on both ends)
$ws.Range("A97").Value = "updel_096"
$ws.Range("B97").Value = "y"
$ws.Range("C97").Value = "两端范围删除"
$ws.Range("D97").Value = "SQLFunction"
$ws.Range("F97").Value = "schema1"
$ws.Range("G97").Value = "updel_value01"
$ws.Range("H97").Value = "delete from `$schema1 where id<5 or id>15"
$ws.Range("I97").Value = "10"
$ws.Range("J97").Value = "select * from `$schema1"
$ws.Range("K97").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_096.csv"
$ws.Range("L97").Value = "csv_containsAll"

# Row 98: updel_097 - 主键not between范围删除
$ws.Range("A98").Value = "updel_097"
$ws.Range("B98").Value = "y"
$ws.Range("C98").Value = "主键not between范围删除"
$ws.Range("D98").Value = "SQLFunction"
$ws.Range("F98").Value = "schema1"
$ws.Range("G98").Value = "updel_value01"
$ws.Range("H98").Value = "delete from `$schema1 where id not between 10 and 11"
$ws.Range("I98").Value = "19"
$ws.Range("J98").Value = "select * from `$schema1"
$ws.Range("K98").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_097.csv"
$ws.Range("L98").Value = "csv_containsAll"

# Row 99: updel_098 - 主键not between范围删除
$ws.Range("A99").Value = "updel_098"
$ws.Range("B99").Value = "y"
$ws.Range("C99").Value = "主键not between范围删除"
$ws.Range("D99").Value = "SQLFunction"
$ws.Range("F99").Value = "schema1"
$ws.Range("G99").Value = "updel_value01"
$ws.Range("H99").Value = "delete from `$schema1 where id not between 10 and 10"
$ws.Range("I99").Value = "20"
$ws.Range("J99").Value = "select * from `$schema1"
$ws.Range("K99").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_098.csv"
$ws.Range("L99").Value = "csv_containsAll"

$ws.Range("K99").Select() | Out-Null
